# --- Edit backend/receipts.xlsx ---------------------------------------
# 1. Rename the existing sheet "2025-04-23" -> "2025-04-24" and update /
#    append its rows of data.
# 2. Add a brand new sheet "2025-04-25" with its own header + 3 rows.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: rename + update data
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "2025-04-24"

# Row 2 updates
$ws1.Range("B2").Value = "Charles Darwin TESTING TWO"
$ws1.Range("C2").Value = "'2025-04-24"
$ws1.Range("C2").Style = "Normal"
$ws1.Range("E2").Value = 13

# Row 3 updates
$ws1.Range("B3").Value = "Justina Wimer"
$ws1.Range("C3").Value = "'2025-04-24"
$ws1.Range("C3").Style = "Normal"
$ws1.Range("E3").Value = 1

# New rows 4-10
$newRows = @(
    @(4, "Matthew Wolz",  "Justina Wimer",      "2025-04-24", "MW", 1, "N/A", "Daily Guest Pass",             3),
    @(5, "James W",       "Greg S",             "2025-04-24", "MW", 2, "N/A", "10 Visit Guest Pass",          25),
    @(6, "N/A",           "Tiffany Neff",       "2025-04-24", "MW", 3, "N/A", "10 Visit Children Guest Pass", 25),
    @(7, "Kafi Rahman",   "Little Kafi Rahman", "2025-04-24", "MW", 4, "N/A", "Youth Guest Pass",             3),
    @(8, "Charles Darwin","Isaac Newton",       "2025-04-24", "MW", 5, "N/A", "Youth Guest Pass",             3),
    @(9, "Matthew Wolz",  "Aiden W",            "2025-04-24", "MW", 6, "N/A", "Daily Guest Pass",             3),
    @(10,"Matthew Wolz",  "Mason Berliner",     "2025-04-24", "MW", 7, "N/A", "Daily Guest Pass",             3)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws1.Range("A$r").Value = $row[1]
    $ws1.Range("B$r").Value = $row[2]
    $ws1.Range("C$r").Value = "'" + $row[3]
    $ws1.Range("C$r").Style = "Normal"
    $ws1.Range("D$r").Value = $row[4]
    $ws1.Range("E$r").Value = $row[5]
    $ws1.Range("F$r").Value = $row[6]
    $ws1.Range("G$r").Value = $row[7]
    $ws1.Range("H$r").Value = $row[8]
}

# ------------------------------------------------------------------
# Sheet 2: brand new sheet, inserted right after sheet 1
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2025-04-25"

$ws2.Range("A1").Value = "Sponsor Name"
$ws2.Range("B1").Value = "Guest Name"
$ws2.Range("C1").Value = "Date"
$ws2.Range("D1").Value = "Initials"
$ws2.Range("E1").Value = "Receipt Number"
$ws2.Range("F1").Value = "Email"
$ws2.Range("G1").Value = "Item"
$ws2.Range("H1").Value = "Price"

$sheet2Rows = @(
    @(2, "matthew wolz", "Adli Jacobs",  "2025-04-25", "MW", 8,  "N/A", "Daily Guest Pass", 3),
    @(3, "Matthew Wolz",  "Aiden Wolz",   "2025-04-25", "MW", 9,  "N/A", "Youth Guest Pass", 3),
    @(4, "matthew wolz",  "jj something", "2025-04-25", "MW", 10, "N/A", "Youth Guest Pass", 3)
)

foreach ($row in $sheet2Rows) {
    $r = $row[0]
    $ws2.Range("A$r").Value = $row[1]
    $ws2.Range("B$r").Value = $row[2]
    $ws2.Range("C$r").Value = "'" + $row[3]
    $ws2.Range("C$r").Style = "Normal"
    $ws2.Range("D$r").Value = $row[4]
    $ws2.Range("E$r").Value = $row[5]
    $ws2.Range("F$r").Value = $row[6]
    $ws2.Range("G$r").Value = $row[7]
    $ws2.Range("H$r").Value = $row[8]
}
